$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The data table tracks lab sample results (lab_id / toc) in columns A:B.
# A block of now-obsolete rows (sample IDs 222236-222242, including the
# "222241 Dup" label row) and a single stray row (sample 222243) are
# removed; all rows below shift up to close the gaps.

# Remove the 8-row block for samples 222236 .. 222242 (rows 57-64)
$ws.Range("A57:B64").EntireRow.Delete() | Out-Null

# Remove the row for sample 222243 (now at row 63 after the shift above)
$ws.Range("A63:B63").EntireRow.Delete() | Out-Null

# Reflect where the user's selection ended up after the edits
$ws.Range("C60").Select() | Out-Null
